# Daily attendance processing - rotate "Recorded By" (column G) entries so
# that the last comma-separated name/email moves to the front of the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $newParts = @($parts[-1]) + $parts[0..($parts.Length-2)]
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value2 = $newVal
        }
    }
}
